$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, reusing the same formatting (style)
# as the other header cells (e.g. G1) by copying formats across.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Mark every data row (2-8) as saved with a value of 1 in column H
$ws.Range("H2:H8").Value = 1
